# edit.ps1
# Applies the "feat: add 2022-Q3 data" change:
#  1. Inserts a new "2022-Q3" sheet right after "总计" (before "2022-Q2"),
#     populated with the Q3 fund-holdings detail table.
#  2. Updates the "总计" (totals) sheet: a new row for 2022-Q3 is inserted
#     at the top of the data (row 2), and every subsequent quarter's row
#     shifts down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: update the "总计" summary sheet with the shifted quarter rows
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$totalsData = @(
  @(2, 0, "2022-Q3", 29, 21.57),
  @(3, 1, "2022-Q2", 60, 10.65),
  @(4, 2, "2022-Q1", 25, 4.21),
  @(5, 3, "2021-Q4", 66, 34.18),
  @(6, 4, "2021-Q3", 110, 67.34),
  @(7, 5, "2021-Q2", 101, 55.6),
  @(8, 6, "2021-Q1", 47, 29.27),
  @(9, 7, "2020-Q4", 124, 53.55),
)

foreach ($row in $totalsData) {
  $r = $row[0]
  $wsTotal.Cells.Item($r, 1).Value = $row[1]
  $wsTotal.Cells.Item($r, 2).Value = $row[2]
  $wsTotal.Cells.Item($r, 3).Value = $row[3]
  $wsTotal.Cells.Item($r, 4).Value = $row[4]
}

# Row 9 is brand-new (the sheet used to stop at row 8) - give column A the
# same bold/bordered style the rest of column A uses before writing its value.
$wsTotal.Cells.Item(8, 1).Copy()
$wsTotal.Cells.Item(9, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(9, 1).Value = 7

# ---------------------------------------------------------------------
# Part 2: add the new "2022-Q3" worksheet (detail table), positioned
# right after "总计" and before "2022-Q2"
# ---------------------------------------------------------------------
$wsQ2Template = $wb.Worksheets.Item("2022-Q2")
$wsQ2Template.Copy($wsQ2Template, $null)
$ws3 = $wb.Worksheets.Item("2022-Q2 (2)")
$ws3.Name = "2022-Q3"

# The template (2022-Q2) has 60 data rows (rows 2-61); 2022-Q3 only needs
# 29 data rows (rows 2-30), drop the extra rows entirely.
$ws3.Range("A31:A61").EntireRow.Delete()

$q3Data = @(
  @(0, "003834", "华夏能源革新股票A", "167.27", "93.80", "5.26", "8.7984", 8),
  @(1, "012967", "广发行业严选三年持有期混合A", "112.21", "94.31", "5.62", "6.3062", 9),
  @(2, "004854", "广发中证全指汽车指数A", "15.34", "94.07", "8.98", "1.3775", 4),
  @(3, "004855", "广发中证全指汽车指数C", "10.15", "94.07", "8.98", "0.9115", 4),
  @(4, "240010", "华宝行业精选混合", "13.17", "88.32", "6.19", "0.8152", 3),
  @(5, "240009", "华宝先进成长混合", "11.36", "87.51", "6.46", "0.7339", 2),
  @(6, "012968", "广发行业严选三年持有期混合C", "10.50", "94.31", "5.62", "0.5901", 9),
  @(7, "516110", "国泰中证800汽车与零部件ETF", "5.83", "97.46", "5.94", "0.3463", 6),
  @(8, "515250", "富国中证智能汽车主题ETF", "6.26", "99.00", "4.99", "0.3124", 6),
  @(9, "013188", "华夏能源革新股票C", "3.47", "93.80", "5.26", "0.1825", 8),
  @(10, "161033", "富国中证智能汽车指数（LOF）", "3.72", "93.93", "4.73", "0.1760", 6),
  @(11, "010956", "天弘中证智能汽车主题指数C", "2.70", "95.20", "4.77", "0.1288", 6),
  @(12, "008056", "南方上证50指数增强A", "1.60", "93.18", "5.68", "0.0909", 6),
  @(13, "002683", "民生加银前沿科技灵活配置混合", "1.50", "92.57", "5.29", "0.0794", 4),
  @(14, "516520", "华泰柏瑞中证智能汽车主题ETF", "1.56", "97.91", "4.88", "0.0761", 6),
  @(15, "013292", "富国中证智能汽车指数(LOF)C", "1.57", "93.93", "4.73", "0.0743", 6),
  @(16, "008057", "南方上证50指数增强C", "1.09", "93.18", "5.68", "0.0619", 6),
  @(17, "007590", "华宝绿色领先股票", "0.90", "89.84", "6.62", "0.0596", 3),
  @(18, "516590", "易方达中证智能电动汽车ETF", "1.39", "98.52", "4.16", "0.0578", 4),
  @(19, "690004", "民生加银稳健成长混合", "1.12", "92.05", "5.10", "0.0571", 4),
  @(20, "516380", "华宝中证智能电动汽车ETF", "1.34", "97.94", "4.13", "0.0553", 4),
  @(21, "159888", "华夏中证智能汽车主题ETF", "1.01", "99.66", "5.02", "0.0507", 6),
  @(22, "002212", "嘉实新起航灵活配置混合", "1.16", "79.78", "4.24", "0.0492", 8),
  @(23, "159889", "国泰中证智能汽车主题ETF", "1.00", "97.21", "4.68", "0.0468", 5),
  @(24, "010955", "天弘中证智能汽车主题指数A", "0.77", "95.20", "4.77", "0.0367", 6),
  @(25, "159795", "汇添富中证智能汽车主题ETF", "0.67", "97.76", "4.94", "0.0331", 6),
  @(26, "159720", "泰康中证智能电动汽车ETF", "0.80", "97.46", "4.13", "0.0330", 4),
  @(27, "159710", "建信中证智能电动汽车ETF", "0.50", "99.27", "4.20", "0.0210", 4),
  @(28, "560000", "浦银安盛中证智能电动汽车ETF", "0.32", "93.36", "4.10", "0.0131", 5),
)

foreach ($row in $q3Data) {
  $r = [int]$row[0] + 2
  $ws3.Cells.Item($r, 1).Value = [int]$row[0]
  $ws3.Cells.Item($r, 2).NumberFormat = "@"
  $ws3.Cells.Item($r, 2).Value = $row[1]
  $ws3.Cells.Item($r, 3).NumberFormat = "@"
  $ws3.Cells.Item($r, 3).Value = $row[2]
  $ws3.Cells.Item($r, 4).NumberFormat = "@"
  $ws3.Cells.Item($r, 4).Value = $row[3]
  $ws3.Cells.Item($r, 5).NumberFormat = "@"
  $ws3.Cells.Item($r, 5).Value = $row[4]
  $ws3.Cells.Item($r, 6).NumberFormat = "@"
  $ws3.Cells.Item($r, 6).Value = $row[5]
  $ws3.Cells.Item($r, 7).NumberFormat = "@"
  $ws3.Cells.Item($r, 7).Value = $row[6]
  $ws3.Cells.Item($r, 8).Value = [int]$row[7]
}
